$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$origStyle = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.229.62'
$ws.Range("D2").Style = $origStyle
$ws.Range("E2").Value = '  +2.30%  '

$origStyle = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.234.17'
$ws.Range("D3").Style = $origStyle
$ws.Range("E3").Value = '  +4.64%  '

$ws.Range("E4").Value = '  -0.01%  '

$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '574.51'
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = '  +1.18%  '

$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.57'
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = '  +7.73%  '

$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = '  -0.17%  '

$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.227.16'
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = '  +4.84%  '

$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.514'
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = '  +3.48%  '

$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.12'
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = '  +9.66%  '

$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.167'
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = '  +5.14%  '

$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.484'
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = '  +3.46%  '

$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.06'
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = '  +5.86%  '

$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000235'
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = '  +2.94%  '

$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.750.78'
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = '  +4.64%  '

$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.210.91'
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = '  +2.22%  '

$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '547.42'
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = '  +10.00%  '

$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.238.19'
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = '  +4.60%  '

$ws.Range("E19").Value = '  +3.01%  '

$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.07'
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = '  +5.54%  '

$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.47'
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = '  +4.63%  '

$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.740'
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = '  +6.61%  '

$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.78'
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = '  +7.73%  '

$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.47'
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = '  +5.98%  '

$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '82.06'
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = '  +3.81%  '

$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = '  -0.08%  '

$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.36'
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = '  +16.11%  '

$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.89'
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = '  +3.73%  '

$ws.Range("E29").Value = '  +8.61%  '

$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '27.88'
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = '  +4.78%  '

$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.75'
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = '  +1.78%  '

$ws.Range("E32").Value = '  +0.21%  '

$ws.Range("E33").Value = '  +4.71%  '

$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '568.49'
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = '  +10.34%  '

$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.77'
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = '  +3.60%  '

$ws.Range("E36").Value = '  +7.65%  '

$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0465'
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = '  +13.71%  '

$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '54.78'
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = '  +2.94%  '

$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0874'
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = '  +9.23%  '

$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.03'
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = '  +11.80%  '

$ws.Range("E41").Value = '  +3.95%  '

$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.134.40'
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = '  +6.31%  '

$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.64'
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = '  +2.52%  '

$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.273'
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = '  +10.38%  '

$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.34'
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = '  +7.72%  '

$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '27.13'
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = '  +6.58%  '

$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₃0566'
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = '  +3.30%  '

$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.999'
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = '  +0.02%  '

$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.113'
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = '  +4.18%  '

$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '122.43'
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = '  +0.61%  '

$ws.Range("B51").Value = 'ThetaToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.25'
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = '  +7.57%  '

